$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated loading_percent results for the 380 kV case
$data = @{
  2 = @{ "B"=18.56572392040826; "C"=9.241994103215211; "E"=11.12051325010195; "F"=16.86991607391233; "G"=38.30126954453143; "H"=16.30901378972903; "L"=9.927261588779871; "N"=17.98528927832327 }
  3 = @{ "B"=18.07574737451615; "C"=8.935776507667962; "E"=11.14971989343208; "F"=15.89584955866808; "G"=38.09353024335096; "H"=16.34251432632526; "L"=9.906247075366476; "N"=18.06039553380951 }
  4 = @{ "B"=17.77276432757274; "C"=8.740554058401969; "E"=11.16918971764735; "F"=15.26997757108491; "G"=37.98157145262777; "H"=16.36721762392421; "L"=9.895290913932671; "N"=18.1085649145385 }
  5 = @{ "B"=17.64895931565053; "C"=8.659255795206828; "E"=11.1775103298204; "F"=15.00819731993403; "G"=37.93989868804854; "H"=16.37831965039403; "L"=9.891318583116135; "N"=18.12871269312987 }
  6 = @{ "B"=17.62838725541474; "C"=8.645653100591897; "E"=11.17891530638038; "F"=14.96433081551593; "G"=37.93321833150809; "H"=16.38022552897414; "L"=9.89068879505766; "N"=18.13208958050876 }
  7 = @{ "B"=17.77109574340553; "C"=8.739464606068397; "E"=11.16930036733337; "F"=15.26647399323137; "G"=37.98099340661291; "H"=16.36736316402484; "L"=9.895235344555427; "N"=18.10883453319285 }
  8 = @{ "B"=18.39734035162626; "C"=9.137943897879056; "E"=11.13026476516313; "F"=16.53996406344768; "G"=38.22642501567222; "H"=16.31970432262192; "L"=9.919613547251489; "N"=18.01076079686172 }
  9 = @{ "B"=19.60016898111747; "C"=9.85957942431963; "E"=11.06590784347489; "F"=19.00274580682531; "G"=38.82973891009819; "H"=16.25922426648451; "L"=9.9827242336016; "N"=17.83464633468364 }
  10 = @{ "B"=20.458182962383; "C"=10.35027608566956; "E"=11.02605430253415; "F"=20.67494806633232; "G"=39.34456671122905; "H"=16.2351212092087; "L"=10.03820710221737; "N"=17.71501021304756 }
  11 = @{ "B"=20.84099248934523; "C"=10.56440977575527; "E"=11.00953653861622; "F"=21.3917225636224; "G"=39.59356933426854; "H"=16.22861273278234; "L"=10.06537298040699; "N"=17.66267600172913 }
  12 = @{ "B"=20.98472848588787; "C"=10.64415339324877; "E"=11.00351343686266; "F"=21.65686569030329; "G"=39.68991697045016; "H"=16.22679173558172; "L"=10.07593179755634; "N"=17.64315675566669 }
  13 = @{ "B"=20.95382907606689; "C"=10.62703950629104; "E"=11.00480030886409; "F"=21.60004134736742; "G"=39.66907657591025; "H"=16.22715525357042; "L"=10.0736457727356; "N"=17.64734732335703 }
  14 = @{ "B"=20.8528430941629; "C"=10.5709974678024; "E"=11.00903636823243; "F"=21.4136618050453; "G"=39.60145515052646; "H"=16.22845000423209; "L"=10.06623624744886; "N"=17.66106416776897 }
  15 = @{ "B"=20.79082251713158; "C"=10.53649402114624; "E"=11.01166126746512; "F"=21.29868154950795; "G"=39.56030052879782; "H"=16.22932697255944; "L"=10.06173291662327; "N"=17.66950496247945 }
  16 = @{ "B"=20.43300162916443; "C"=10.33609577928721; "E"=11.0271662145511; "F"=20.62722412089977; "G"=39.32858558785541; "H"=16.23563647321939; "L"=10.036470065881; "N"=17.71847225665462 }
  17 = @{ "B"=20.21145962119371; "C"=10.21080236535321; "E"=11.03709085074551; "F"=20.20408069617459; "G"=39.19017758627812; "H"=16.24065060639319; "L"=10.02146185875633; "N"=17.74904583469467 }
  18 = @{ "B"=20.08333602553486; "C"=10.13788353917894; "E"=11.04295094815288; "F"=19.95656407809808; "G"=39.11196722673409; "H"=16.24395388904696; "L"=10.0130112041696; "N"=17.76682766808658 }
  19 = @{ "B"=20.03984027906397; "C"=10.11304913474194; "E"=11.04496113134473; "F"=19.87204792380562; "G"=39.08572883335586; "H"=16.24514424384398; "L"=10.01018131591641; "N"=17.77288213290451 }
  20 = @{ "B"=20.23511653816779; "C"=10.2242286706826; "E"=11.03601865467198; "F"=20.24955283636157; "G"=39.20476715231966; "H"=16.24007342726227; "L"=10.02304074341476; "N"=17.74577087845965 }
  21 = @{ "B"=20.88253949337135; "C"=10.5874951280806; "E"=11.00778584316116; "F"=21.46857628470577; "G"=39.62126202006974; "H"=16.22805221616497; "L"=10.06840527539612; "N"=17.65702710880583 }
  22 = @{ "B"=21.29846465108217; "C"=10.81706187846464; "E"=10.99068529178872; "F"=22.22866616901552; "G"=39.90541032640508; "H"=16.22394822180894; "L"=10.09963450017115; "N"=17.60076747695755 }
  23 = @{ "B"=21.0771817781087; "C"=10.69526711406234; "E"=10.99968852837591; "F"=21.82633154475857; "G"=39.75268762092865; "H"=16.22579440358724; "L"=10.08282407994629; "N"=17.63063572142854 }
  24 = @{ "B"=20.22442358994676; "C"=10.21816139482146; "E"=11.03650291392838; "F"=20.22900810905287; "G"=39.19816696077694; "H"=16.24033306004221; "L"=10.02232637555863; "N"=17.74725084846888 }
  25 = @{ "B"=19.27858516687954; "C"=9.671096408452275; "E"=11.08201318036395; "F"=18.34778573295691; "G"=38.65372507667592; "H"=16.27203038686702; "L"=9.96403317525937; "N"=17.88056782122947 }
}

foreach ($row in $data.Keys) {
  $rowData = $data[$row]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$row").Value = $rowData[$col]
  }
}
